$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.032.11"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.825.89"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4394"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07269"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8441"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "1.785.76"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.652"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07067"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.298"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008783"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "27.061.47"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.154"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "2.052.19"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.982"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.208"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08801"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7408"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.905"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.427"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05238"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.241"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.866"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5158"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1696"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.533"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4808"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.930"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06335"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.658"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
